$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) column before column C, shifting the old
# "Bet amount..Net Result" block from C:H to D:I. Column A (Date) and
# column B (Leg info text) stay put.
$ws.Columns("C:C").Insert()

# New column B becomes "Player" (mostly blank below), old B1 header
# text ("Leg info") moves into the freshly inserted C1.
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "Leg info"

# --- New rows of bet data (rows 10-18) ---

# Row 10 - Kawhi, still has a date like the rows above it
$ws.Range("A10").Value = 45382
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("B10").Value = "Kawhi 20+ pts"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1.14
$ws.Range("F10").Value = 1.4
$ws.Range("G10").Formula = "=F10/E10-1"

# Row 11 - Austin Reeves (no date)
$ws.Range("B11").Value = "Austin Reeves 15+ pts"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1.48
$ws.Range("F11").Value = 1.74
$ws.Range("G11").Formula = "=F11/E11-1"

# Row 12 - Rui Hachimura (no date)
$ws.Range("B12").Value = "Rui Hachimura 4+ rebs"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1.5
$ws.Range("F12").Value = 1.71
$ws.Range("G12").Formula = "=F12/E12-1"

# Row 13 - Lebron (no date)
$ws.Range("B13").Value = "Lebron 20+ pts"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1.11
$ws.Range("F13").Value = 1.3
$ws.Range("G13").Formula = "=F13/E13-1"

# Row 14 - D'angelo Russel (no date)
$ws.Range("B14").Value = "D'angelo Russel 6+ assits"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1.42
$ws.Range("F14").Value = 2.1
$ws.Range("G14").Formula = "=F14/E14-1"

# Row 15 - Nic Claxton 10+ pts (no date, uses the new "Leg info" column C)
$ws.Range("B15").Value = "Nic Claxton"
$ws.Range("C15").Value = "10+ pts"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1.34
$ws.Range("F15").Value = 1.51
$ws.Range("G15").Formula = "=F15/E15-1"

# Row 16 - Mikal Bridges 15+ pts
$ws.Range("B16").Value = "Mikal Bridges"
$ws.Range("C16").Value = "15+ pts"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1.04
$ws.Range("F16").Value = 1.32
$ws.Range("G16").Formula = "=F16/E16-1"

# Row 17 - Mikal Bridges 6+ rebs
$ws.Range("B17").Value = "Mikal Bridges"
$ws.Range("C17").Value = "6+ rebs"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 2.57
$ws.Range("F17").Value = 3.3
$ws.Range("G17").Formula = "=F17/E17-1"

# Row 18 - Terry Rozier 15+ pts
$ws.Range("B18").Value = "Terry Rozier"
$ws.Range("C18").Value = "15+ pts"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1.06
$ws.Range("F18").Value = 1.56
$ws.Range("G18").Formula = "=F18/E18-1"

$ws.Range("B1").Select()
